$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$abstract104 = @"
Novel coronavirus (SARS-CoV-2) causing COVID-19 disease has arisen to be a pandemic.
 Since there is a close association between other viral infection cases by epidemics and environmental factors, this study intends to unveil meteorological effects on the outbreak of COVID-19 across eight divisions of Bangladesh from March to April 2020. A compound Poisson generalized linear modeling (CPGLM), along with a Monte-Carlo method and random forest (RF) model, was employed to explore how meteorological factors affecting the COVID-19 transmission in Bangladesh.
 Results showed that subtropical climate (mean temperature about 26.6 °C, mean relative humidity (MRH) 64%, and rainfall approximately 3 mm) enhanced COVD-19 onset.
 The CPGLM model revealed that every 1 mm increase in rainfall elevated by 30.99% (95% CI 77.18%, − 15.20%) COVID-19 cases, while an increase of 1 °C of diurnal temperature (TDN) declined the confirmed cases by − 14.2% (95% CI 9.73%, − 38.13%) on the lag 1 and lag 2, respectively.
 In addition, NRH and MRH had the highest increase (17.98% (95% CI 22.5%, 13.42%) and 19.92% (95% CI: 25.71%, 14.13%)) of COVID-19 cased in lag 4. The results of the RF model indicated that TDN and AH (absolute humidity) influence the COVID-19 cases most.
 In the Dhaka division, MRH is the most vital meteorological factor that affects COVID-19 deaths.
 This study indicates the humidity and rainfall are crucial factors affecting the COVID-19 case, which is contrary to many previous studies in other countries.
 These outcomes can have policy formulation for the suppression of the COVID-19 outbreak in Bangladesh.

"@

$authors105 = "[Abu Reza Md. Towfiqul%Islam%towfiq_dm@brur.ac.bd%1,    Md.%Hasanuzzaman%NULL%2,    Md.%Hasanuzzaman%NULL%0,    Mashura%Shammi%NULL%1,    Roquia%Salam%NULL%1,    Md.%Bodrud-Doza%NULL%1,    Md. Mostafizur%Rahman%rahmanmm@juniv.edu%1,    Md. Abdul%Mannan%NULL%1,    Saleemul%Huq%NULL%1]"

$abstract106 = @"
The emergence of coronavirus disease 2019 (COVID-19) has become a worldwide pandemic after its first outbreak in Wuhan, China.
 However, it remains unclear whether COVID-19 death is linked to ambient air pollutants or meteorological conditions.
 We collected the daily COVID-19 death number, air quality index (AQI), ambient air pollutant concentrations, and meteorological variables data of Wuhan between Jan 25 and April 7, 2020. The Pearson and Poisson regression models were used accordingly to understand the association between COVID-19 deaths and each risk factor.
 The daily COVID-19 deaths were positively correlated with AQI (slope = 0.4 ± 0.09, R2 = 0.24, p &lt; 0.01).
 Detailedly, PM2.5 was the only pollutant exhibiting a positive association (relative risk (RR) = 1.079, 95%CI 1.071–1.086, p &lt; 0.01) with COVID-19 deaths.
 The PM10, SO2, and CO were all also significantly associated with COVID-19 deaths, but in negative pattern (p &lt; 0.01).
 Among them, PM10 and CO had the highest and lowest RR, which equaled to 0.952 (95%CI 0.945–0.959) and 0.177 (95%CI 0.131–0.24), respectively.
 Additionally, temperature was inversely associated with COVID-19 deaths (RR = 0.861, 95%CI 0.851–0.872, p &lt; 0.01).
 Contrarily, diurnal temperature range was positively associated with COVID-19 deaths (RR = 1.014, 95%CI 1.003–1.025, p &lt; 0.05).
 The data suggested that PM2.5 and diurnal temperature range are tightly associated with increased COVID-19 deaths.

"@

$authors107 = "[Ying%Jiang%yjiang1@vt.edu%1,    Jing%Xu%xujingjingalpha@sohu.com%0]"

$ws.Range("D8").Value = $abstract104
$ws.Range("E8").Value = $authors105
$ws.Range("D9").Value = $abstract106
$ws.Range("E9").Value = $authors107
